$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the total opex values in D2:AH4 with the new sensitivity-ratio figures
# derived from the exact 1-year shifting (no double-deployment capex) IEV model run.

# Row 2
$ws.Range("D2").Value = 8275.84085119876
$ws.Range("E2").Value = 8181.1364275615206
$ws.Range("F2").Value = 8089.7048406415706
$ws.Range("G2").Value = 7771.0890741038456
$ws.Range("H2").Value = 7794.1057507254573
$ws.Range("I2").Value = 7682.3896806233615
$ws.Range("J2").Value = 7694.0204871794995
$ws.Range("K2").Value = 7731.8609711502249
$ws.Range("L2").Value = 7703.2236325364447
$ws.Range("M2").Value = 7747.5870479553641
$ws.Range("N2").Value = 7629.195061081562
$ws.Range("O2").Value = 7603.4961444462206
$ws.Range("P2").Value = 7594.1035764397275
$ws.Range("Q2").Value = 7602.2299267878925
$ws.Range("R2").Value = 7629.2420571742423
$ws.Range("S2").Value = 7676.6777567949457
$ws.Range("T2").Value = 7745.906425448622
$ws.Range("U2").Value = 7838.9837546686895
$ws.Range("V2").Value = 7958.0245797795915
$ws.Range("W2").Value = 8105.3928308304221
$ws.Range("X2").Value = 8283.7714983540518
$ws.Range("Y2").Value = 8496.0564634173988
$ws.Range("Z2").Value = 8745.5222803205561
$ws.Range("AA2").Value = 9035.8062146272969
$ws.Range("AB2").Value = 9370.9389121799286
$ws.Range("AC2").Value = 9755.3740418586221
$ws.Range("AD2").Value = 10194.015030355422
$ws.Range("AE2").Value = 10692.236116861128
$ws.Range("AF2").Value = 11255.893697186901
$ws.Range("AG2").Value = 11891.322151846314
$ws.Range("AH2").Value = 12449.998927407036

# Row 3
$ws.Range("D3").Value = 8276.2262983212659
$ws.Range("E3").Value = 8182.3162547510738
$ws.Range("F3").Value = 8092.303188744987
$ws.Range("G3").Value = 7775.6813845101387
$ws.Range("H3").Value = 7801.6248867715558
$ws.Range("I3").Value = 7691.219710630493
$ws.Range("J3").Value = 7706.3994978338542
$ws.Range("K3").Value = 7748.7815048128496
$ws.Range("L3").Value = 7725.6331137970919
$ws.Range("M3").Value = 7776.3881995830397
$ws.Range("N3").Value = 7708.5095826346669
$ws.Range("O3").Value = 7739.671792174403
$ws.Range("P3").Value = 7790.8964786619745
$ws.Range("Q3").Value = 7863.4651268095104
$ws.Range("R3").Value = 7958.466917087776
$ws.Range("S3").Value = 8077.4019300108184
$ws.Range("T3").Value = 8221.4249014095585
$ws.Range("U3").Value = 8392.1787011638553
$ws.Range("V3").Value = 8591.2407388316096
$ws.Range("W3").Value = 8820.5800250535922
$ws.Range("X3").Value = 9081.9005902509052
$ws.Range("Y3").Value = 9377.1442375225906
$ws.Range("Z3").Value = 9708.2320929987054
$ws.Range("AA3").Value = 10077.159442694563
$ws.Range("AB3").Value = 10485.63115065762
$ws.Range("AC3").Value = 10935.373911165034
$ws.Range("AD3").Value = 11427.718439416813
$ws.Range("AE3").Value = 11963.642696498273
$ws.Range("AF3").Value = 12543.364621557872
$ws.Range("AG3").Value = 13166.26940484773
$ws.Range("AH3").Value = 13567.277373437606

# Row 4
$ws.Range("D4").Value = 8276.9916239738759
$ws.Range("E4").Value = 8185.0916610053919
$ws.Range("F4").Value = 8098.3696655893173
$ws.Range("G4").Value = 7786.8038883910813
$ws.Range("H4").Value = 7819.4669257335254
$ws.Range("I4").Value = 7712.2285314956071
$ws.Range("J4").Value = 7735.8127894905883
$ws.Range("K4").Value = 7789.1284401704752
$ws.Range("L4").Value = 7779.0353661122526
$ws.Range("M4").Value = 7844.9609335224359
$ws.Range("N4").Value = 7824.4965856560702
$ws.Range("O4").Value = 7915.3634065639508
$ws.Range("P4").Value = 8032.6638878805297
$ws.Range("Q4").Value = 8177.8376849314309
$ws.Range("R4").Value = 8352.1058344555913
$ws.Range("S4").Value = 8556.858165963853
$ws.Range("T4").Value = 8793.1761510267388
$ws.Range("U4").Value = 9062.738969207232
$ws.Range("V4").Value = 9366.8342241571263
$ws.Range("W4").Value = 9706.9687987724101
$ws.Range("X4").Value = 10084.467269757342
$ws.Range("Y4").Value = 10500.057576060806
$ws.Range("Z4").Value = 10954.652487216095
$ws.Range("AA4").Value = 11448.307827869448
$ws.Range("AB4").Value = 11980.398814804015
$ws.Range("AC4").Value = 12549.230148925035
$ws.Range("AD4").Value = 13151.661769889599
$ws.Range("AE4").Value = 13782.627326299455
$ws.Range("AF4").Value = 14434.716233262759
$ws.Range("AG4").Value = 15096.982303227551
$ws.Range("AH4").Value = 15294.97702322077

Write-Host "Updated total opex sensitivity values in D2:AH4"
